$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Hunk 1+2: paragraph "php artisan make:middleware Name" -------------
# Remove the rFonts w:ascii="Arial"/w:hAnsi="Arial" overrides from both
# runs of this paragraph (tab run + text run).
$rng1 = $d.Content
$found = $rng1.Find.Execute("php artisan make:middleware Name", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "paragraph 'make:middleware Name' not found" }
$p1 = $rng1.Paragraphs(1)
$xml1 = "<w:p $wns>" +
          "<w:pPr>" +
            "<w:pStyle w:val=`"Normal`"/>" +
            "<w:spacing w:lineRule=`"auto`" w:line=`"360`" w:before=`"0`" w:after=`"0`"/>" +
            "<w:jc w:val=`"left`"/>" +
            "<w:rPr></w:rPr>" +
          "</w:pPr>" +
          "<w:r>" +
            "<w:rPr>" +
              "<w:b w:val=`"false`"/>" +
              "<w:bCs w:val=`"false`"/>" +
              "<w:u w:val=`"none`"/>" +
              "<w:lang w:val=`"zxx`" w:eastAsia=`"zxx`" w:bidi=`"zxx`"/>" +
            "</w:rPr>" +
            "<w:tab/>" +
          "</w:r>" +
          "<w:r>" +
            "<w:rPr>" +
              "<w:rStyle w:val=`"Style14`"/>" +
              "<w:rFonts w:eastAsia=`"Calibri`" w:cs=`"`" w:cstheme=`"minorBidi`" w:eastAsiaTheme=`"minorHAnsi`"/>" +
              "<w:b w:val=`"false`"/>" +
              "<w:bCs w:val=`"false`"/>" +
              "<w:color w:val=`"auto`"/>" +
              "<w:kern w:val=`"0`"/>" +
              "<w:sz w:val=`"22`"/>" +
              "<w:szCs w:val=`"22`"/>" +
              "<w:u w:val=`"none`"/>" +
              "<w:lang w:val=`"zxx`" w:eastAsia=`"zxx`" w:bidi=`"zxx`"/>" +
            "</w:rPr>" +
            "<w:t>php artisan make:middleware Name</w:t>" +
          "</w:r>" +
        "</w:p>"
$p1.Range.InsertXML($xml1)

# --- Hunk 3: paragraph "php artisan make:request Name" ------------------
# Merge the tab-run and the text-run into a single run.
$rng2 = $d.Content
$found = $rng2.Find.Execute("php artisan make:request Name", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "paragraph 'make:request Name' not found" }
$p2 = $rng2.Paragraphs(1)
$xml2 = "<w:p $wns>" +
          "<w:pPr>" +
            "<w:pStyle w:val=`"Normal`"/>" +
            "<w:spacing w:lineRule=`"auto`" w:line=`"360`" w:before=`"0`" w:after=`"0`"/>" +
            "<w:jc w:val=`"left`"/>" +
            "<w:rPr></w:rPr>" +
          "</w:pPr>" +
          "<w:r>" +
            "<w:rPr>" +
              "<w:rStyle w:val=`"Style14`"/>" +
              "<w:rFonts w:eastAsia=`"Calibri`" w:cs=`"`" w:cstheme=`"minorBidi`" w:eastAsiaTheme=`"minorHAnsi`"/>" +
              "<w:b w:val=`"false`"/>" +
              "<w:bCs w:val=`"false`"/>" +
              "<w:color w:val=`"auto`"/>" +
              "<w:kern w:val=`"0`"/>" +
              "<w:sz w:val=`"22`"/>" +
              "<w:szCs w:val=`"22`"/>" +
              "<w:u w:val=`"none`"/>" +
              "<w:lang w:val=`"zxx`" w:eastAsia=`"zxx`" w:bidi=`"zxx`"/>" +
            "</w:rPr>" +
            "<w:tab/>" +
            "<w:t>php artisan make:request Name</w:t>" +
          "</w:r>" +
        "</w:p>"
$p2.Range.InsertXML($xml2)

# --- Hunk 4: new paragraph after "php artisan lang:publish" -------------
$rng3 = $d.Content
$found = $rng3.Find.Execute("php artisan lang:publish", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "paragraph 'lang:publish' not found" }
$p3 = $rng3.Paragraphs(1)
$endRng = $p3.Range
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$xml3 = "<w:p $wns>" +
          "<w:pPr>" +
            "<w:pStyle w:val=`"Normal`"/>" +
            "<w:spacing w:lineRule=`"auto`" w:line=`"360`" w:before=`"0`" w:after=`"0`"/>" +
            "<w:jc w:val=`"left`"/>" +
            "<w:rPr>" +
              "<w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`"/>" +
            "</w:rPr>" +
          "</w:pPr>" +
          "<w:r>" +
            "<w:rPr>" +
              "<w:b w:val=`"false`"/>" +
              "<w:bCs w:val=`"false`"/>" +
              "<w:u w:val=`"none`"/>" +
            "</w:rPr>" +
            "<w:tab/>" +
          "</w:r>" +
          "<w:r>" +
            "<w:rPr>" +
              "<w:b w:val=`"false`"/>" +
              "<w:bCs w:val=`"false`"/>" +
              "<w:u w:val=`"none`"/>" +
            "</w:rPr>" +
            "<w:t>https://laravel.com/docs/10.x/localization</w:t>" +
          "</w:r>" +
        "</w:p>"
$lastPara.Range.InsertXML($xml3)

Write-Host "done"
